$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update volume number and date range strings (new reporting week)
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# Fix C14: was mistakenly stored as the text "0"; now a real number
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"

# Row 14: Murder
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -66.666666666666
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -63.636363636363
$ws.Range("I14").Value = 34
$ws.Range("J14").Value = 43
$ws.Range("K14").Value = -20.930232558139
$ws.Range("L14").Value = -29.166666666666
$ws.Range("M14").Value = -54.054054054054
$ws.Range("N14").Value = -86.178861788617

# Row 15: Rape
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = -62.5
$ws.Range("F15").Value = 19
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = -5
$ws.Range("I15").Value = 124
$ws.Range("J15").Value = 126
$ws.Range("K15").Value = -1.587301587301
$ws.Range("L15").Value = 0.813008130081
$ws.Range("M15").Value = 10.714285714285
$ws.Range("N15").Value = -63.421828908554

# Row 16: Robbery
$ws.Range("C16").Value = 50
$ws.Range("D16").Value = 58
$ws.Range("E16").Value = -13.793103448275
$ws.Range("F16").Value = 179
$ws.Range("G16").Value = 228
$ws.Range("H16").Value = -21.491228070175
$ws.Range("I16").Value = 1270
$ws.Range("J16").Value = 1386
$ws.Range("K16").Value = -8.369408369408
$ws.Range("L16").Value = 22.23291626564
$ws.Range("M16").Value = -30.143014301430
$ws.Range("N16").Value = -85.333179350964

# Row 17: Fel. Assault
$ws.Range("C17").Value = 89
$ws.Range("D17").Value = 84
$ws.Range("E17").Value = 5.952380952380
$ws.Range("F17").Value = 367
$ws.Range("G17").Value = 375
$ws.Range("H17").Value = -2.133333333333
$ws.Range("I17").Value = 2284
$ws.Range("J17").Value = 2200
$ws.Range("K17").Value = 3.818181818181
$ws.Range("L17").Value = 27.384272169548
$ws.Range("M17").Value = 27.171492204899
$ws.Range("N17").Value = -49.591701611123

# Row 18: Burglary
$ws.Range("C18").Value = 44
$ws.Range("D18").Value = 46
$ws.Range("E18").Value = -4.347826086956
$ws.Range("F18").Value = 156
$ws.Range("G18").Value = 161
$ws.Range("H18").Value = -3.105590062111
$ws.Range("I18").Value = 1083
$ws.Range("J18").Value = 1289
$ws.Range("K18").Value = -15.981380915438
$ws.Range("L18").Value = 10.397553516819
$ws.Range("M18").Value = -30.353697749196
$ws.Range("N18").Value = -82.934131736527

# Row 19: Gr. Larceny
$ws.Range("C19").Value = 117
$ws.Range("D19").Value = 136
$ws.Range("E19").Value = -13.970588235294
$ws.Range("F19").Value = 478
$ws.Range("G19").Value = 531
$ws.Range("H19").Value = -9.981167608286
$ws.Range("I19").Value = 3097
$ws.Range("J19").Value = 3096
$ws.Range("K19").Value = 0.032299741602
$ws.Range("L19").Value = 35.773783428320
$ws.Range("M19").Value = 45.947219604147
$ws.Range("N19").Value = -10.672050764349

# Row 20: G.L.A.
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 43
$ws.Range("E20").Value = 6.976744186046
$ws.Range("F20").Value = 152
$ws.Range("G20").Value = 148
$ws.Range("H20").Value = 2.702702702702
$ws.Range("I20").Value = 942
$ws.Range("J20").Value = 935
$ws.Range("K20").Value = 0.748663101604
$ws.Range("L20").Value = 25.935828877005
$ws.Range("M20").Value = 23.298429319371
$ws.Range("N20").Value = -81.456692913385

# Row 21: TOTAL
$ws.Range("C21").Value = 350
$ws.Range("D21").Value = 378
$ws.Range("E21").Value = -7.407407407407
$ws.Range("F21").Value = 1355
$ws.Range("G21").Value = 1474
$ws.Range("H21").Value = -8.073270013568
$ws.Range("I21").Value = 8834
$ws.Range("J21").Value = 9075
$ws.Range("K21").Value = -2.655647382920
$ws.Range("L21").Value = 25.966063025809
$ws.Range("M21").Value = 7.195728673704
$ws.Range("N21").Value = -69.185154178875

# Row 22: Transit
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = -27.272727272727
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 31
$ws.Range("H22").Value = -16.129032258064
$ws.Range("I22").Value = 158
$ws.Range("J22").Value = 203
$ws.Range("K22").Value = -22.167487684729
$ws.Range("L22").Value = 26.4
$ws.Range("M22").Value = -28.828828828828
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 33
$ws.Range("E23").Value = -3.030303030303
$ws.Range("F23").Value = 123
$ws.Range("G23").Value = 133
$ws.Range("H23").Value = -7.518796992481
$ws.Range("I23").Value = 873
$ws.Range("J23").Value = 811
$ws.Range("K23").Value = 7.644882860665
$ws.Range("L23").Value = 17.496635262449
$ws.Range("M23").Value = 46.969696969697
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 245
$ws.Range("D24").Value = 263
$ws.Range("E24").Value = -6.844106463878
$ws.Range("F24").Value = 1044
$ws.Range("G24").Value = 1074
$ws.Range("H24").Value = -2.793296089385
$ws.Range("I24").Value = 6655
$ws.Range("J24").Value = 6967
$ws.Range("K24").Value = -4.478254628965
$ws.Range("L24").Value = 25.282379518072
$ws.Range("M24").Value = 24.253174010455
$ws.Range("N24").Value = "***.*"

# Row 25: Misd. Assault
$ws.Range("C25").Value = 122
$ws.Range("D25").Value = 96
$ws.Range("E25").Value = 27.083333333333
$ws.Range("F25").Value = 546
$ws.Range("G25").Value = 447
$ws.Range("H25").Value = 22.147651006711
$ws.Range("I25").Value = 3315
$ws.Range("J25").Value = 3302
$ws.Range("K25").Value = 0.393700787401
$ws.Range("L25").Value = 40.406607369758
$ws.Range("M25").Value = -23.174971031286
$ws.Range("N25").Value = "***.*"

# Row 26: UCR Rape*
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 23.076923076923
$ws.Range("I26").Value = 188
$ws.Range("J26").Value = 197
$ws.Range("K26").Value = -4.568527918781
$ws.Range("L26").Value = -9.178743961352
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = -44.444444444444
$ws.Range("F27").Value = 45
$ws.Range("G27").Value = 60
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 330
$ws.Range("J27").Value = 339
$ws.Range("K27").Value = -2.654867256637
$ws.Range("L27").Value = -9.836065573770
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = -37.5
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 46
$ws.Range("H28").Value = -56.521739130434
$ws.Range("I28").Value = 126
$ws.Range("J28").Value = 180
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -38.536585365853
$ws.Range("M28").Value = -53.505535055350
$ws.Range("N28").Value = -87.743190661478

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 38
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 110
$ws.Range("J29").Value = 152
$ws.Range("K29").Value = -27.631578947368
$ws.Range("L29").Value = -36.416184971098
$ws.Range("M29").Value = -48.356807511737
$ws.Range("N29").Value = -88.133764832794

# Row 30: Hate Crimes
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 20
$ws.Range("I30").Value = 37
$ws.Range("J30").Value = 42
$ws.Range("K30").Value = -11.904761904761
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"
